$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# C1: Main PID -> Main GUPRI (i.e. PID)
$ws.Range("C1").Value = "Main GUPRI (i.e. PID) under which all terms are defined. Preference is to use PURLs or W3IDs as they provide permanent resolvable identifiers."

# D4: add "is" after Onotlogy
$ws.Range("D4").Value = "Prefix for SKOS Onotlogy is which our base for defining SKOS based controlled vocabulary"

# C5: add trailing space
$ws.Range("C5").Value = "http://purl.org/pav/ "

# C6: add trailing space
$ws.Range("C6").Value = "http://purl.org/dc/terms/ "

# D8: state -> define, unit -> units
$ws.Range("D8").Value = "Prefix for QUDT Vocabulary of Units which terms we will use to semantically define units of measurements for terms (i.e., variables) we are defining "

# I19: clear the cell content
$ws.Range("I19").Value = ""

# H20: update URL
$ws.Range("H20").Value = "https://mmisw.org/ont/cf/parameter/wind_speed"
